# "Generate Report for Handoff"
#
# A new handoff was just produced for the 0ee8737e-520c-4d47-aee1-f86d5f8d8138
# source file (row 4 of the status tables), so its "Latest Handoff Datetime"
# column (D) is refreshed to the new handoff timestamp on both the zh-cn and
# de-de status sheets. (Row 5 - the 3b9709a7... file - keeps its previous
# handoff timestamp unchanged.)

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-18 03:19:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-18 03:19:49"
